$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 12134.4
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 12134.4
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 12134.4
$ws.Range("M86").ClearContents() | Out-Null
$ws.Range("N86").Value = -14380.4
$ws.Range("H89").Value = 12134.4
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 12134.4
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 60672
$ws.Range("M89").ClearContents() | Out-Null
$ws.Range("N89").Value = -71904
$ws.Range("H105").Value = 33890.332
$ws.Range("J105").Value = 33890.332
$ws.Range("L105").Value = 33890.332
$ws.Range("N105").Value = -40878.332
$ws.Range("H111").Value = 2215.4546
$ws.Range("I111").Value = 1371.1666
$ws.Range("J111").Value = 3228.6
$ws.Range("K111").Value = 4113.4998
$ws.Range("L111").Value = 9685.799999999999
$ws.Range("M111").Value = -1046.4998
$ws.Range("N111").Value = -15819.8
$ws.Range("H137").Value = 1226.289
$ws.Range("I137").Value = 1016.2619
$ws.Range("J137").Value = 4166.6665
$ws.Range("K137").Value = 3048.7857
$ws.Range("L137").Value = 12499.9995
$ws.Range("M137").Value = -498.7856999999999
$ws.Range("N137").Value = -17599.9995
$ws.Range("H138").Value = 5931.971
$ws.Range("I138").Value = 1591.5217
$ws.Range("K138").Value = 4774.5651
$ws.Range("M138").Value = 365.4349000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 90
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents() | Out-Null
$ws.Range("H32").Value = 19400.643
$ws.Range("I32").Value = 3441.5308
$ws.Range("K32").Value = 3441.5308
$ws.Range("M32").Value = -3154.5308
$ws.Range("H45").Value = 46781.684
$ws.Range("I45").Value = 100958.2
$ws.Range("J45").Value = 1634.5834
$ws.Range("K45").Value = 100958.2
$ws.Range("L45").Value = 1634.5834
$ws.Range("M45").Value = -100581.2
$ws.Range("N45").Value = -2388.5834
$ws.Range("H61").Value = 2146.205
$ws.Range("I61").Value = 1460.6316
$ws.Range("J61").Value = 2797.5
$ws.Range("K61").Value = 1460.6316
$ws.Range("L61").Value = 2797.5
$ws.Range("M61").Value = -1248.6316
$ws.Range("N61").Value = -3221.5
$ws.Range("H88").Value = 3077.7778
$ws.Range("I88").Value = 3400
$ws.Range("J88").Value = 2985.7144
$ws.Range("K88").Value = 3400
$ws.Range("L88").Value = 2985.7144
$ws.Range("M88").Value = -2994
$ws.Range("N88").Value = -3797.7144
$ws.Range("H91").Value = 3077.7778
$ws.Range("I91").Value = 3400
$ws.Range("J91").Value = 2985.7144
$ws.Range("K91").Value = 3400
$ws.Range("L91").Value = 2985.7144
$ws.Range("M91").Value = -1996
$ws.Range("N91").Value = -5793.7144
$ws.Range("H132").Value = 2266.58
$ws.Range("I132").Value = 2365.2827
$ws.Range("J132").Value = 1131.5
$ws.Range("K132").Value = 7095.848100000001
$ws.Range("L132").Value = 3394.5
$ws.Range("M132").Value = -4565.848100000001
$ws.Range("N132").Value = -8454.5
$ws.Range("H136").Value = 2146.205
$ws.Range("I136").Value = 1460.6316
$ws.Range("J136").Value = 2797.5
$ws.Range("K136").Value = 4381.8948
$ws.Range("L136").Value = 8392.5
$ws.Range("M136").Value = -1831.8948
$ws.Range("N136").Value = -13492.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 90
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents() | Out-Null
$ws.Range("H26").Value = 20249.5
$ws.Range("I26").Value = 20249.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 20249.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -19957.5
$ws.Range("N26").ClearContents() | Out-Null
$ws.Range("H86").Value = 43790.703
$ws.Range("I86").Value = 67707.17999999999
$ws.Range("J86").Value = 3132.7
$ws.Range("K86").Value = 67707.17999999999
$ws.Range("L86").Value = 3132.7
$ws.Range("M86").Value = -66584.17999999999
$ws.Range("N86").Value = -5378.7
$ws.Range("H89").Value = 43790.703
$ws.Range("I89").Value = 67707.17999999999
$ws.Range("J89").Value = 3132.7
$ws.Range("K89").Value = 338535.9
$ws.Range("L89").Value = 15663.5
$ws.Range("M89").Value = -332919.9
$ws.Range("N89").Value = -26895.5
$ws.Range("H134").Value = 2270.976
$ws.Range("I134").Value = 2278.639
$ws.Range("J134").Value = 2225
$ws.Range("K134").Value = 6835.917
$ws.Range("L134").Value = 6675
$ws.Range("M134").Value = -4300.917
$ws.Range("N134").Value = -11745

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 16000
$ws.Range("J57").Value = 16000
$ws.Range("L57").Value = 16000
$ws.Range("N57").Value = -17120
$ws.Range("H62").Value = 2349.125
$ws.Range("J62").Value = 2349.125
$ws.Range("L62").Value = 2349.125
$ws.Range("N62").Value = -3597.125
$ws.Range("H65").Value = 2349.125
$ws.Range("J65").Value = 2349.125
$ws.Range("L65").Value = 11745.625
$ws.Range("N65").Value = -17985.625
$ws.Range("H88").Value = 43480
$ws.Range("J88").Value = 43480
$ws.Range("L88").Value = 43480
$ws.Range("N88").Value = -44292
$ws.Range("H91").Value = 43480
$ws.Range("J91").Value = 43480
$ws.Range("L91").Value = 43480
$ws.Range("N91").Value = -46288
$ws.Range("H132").Value = 4515.6665
$ws.Range("I132").Value = 4567.857
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 13703.571
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -11173.571
$ws.Range("N132").Value = -18059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100000296
$ws.Range("I4").Value = 118.75
$ws.Range("K4").Value = 356.25
$ws.Range("M4").Value = -244.25
$ws.Range("H107").Value = 1112.0769
$ws.Range("I107").Value = 510
$ws.Range("K107").Value = 1530
$ws.Range("M107").Value = 390
$ws.Range("H120").Value = 337310
$ws.Range("I120").Value = 337310
$ws.Range("K120").Value = 1011930
$ws.Range("M120").Value = -1007092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 40765
$ws.Range("J68").Value = 40765
$ws.Range("L68").Value = 40765
$ws.Range("N68").Value = -42387
$ws.Range("H71").Value = 40765
$ws.Range("J71").Value = 40765
$ws.Range("L71").Value = 122295
$ws.Range("N71").Value = -130407
$ws.Range("H80").Value = 200207400
$ws.Range("I80").Value = 250258750
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 250258750
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -250257752
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 200207400
$ws.Range("I83").Value = 250258750
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 1251293750
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -1251288758
$ws.Range("N83").Value = -19984
$ws.Range("H107").Value = 594914.75
$ws.Range("I107").Value = 467
$ws.Range("J107").Value = 1444125.9
$ws.Range("K107").Value = 467
$ws.Range("L107").Value = 1444125.9
$ws.Range("M107").Value = 1453
$ws.Range("N107").Value = -1447965.9
$ws.Range("H113").Value = 1736.8
$ws.Range("I113").Value = 1704
$ws.Range("K113").Value = 1704
$ws.Range("M113").Value = 466
$ws.Range("H132").Value = 2263.9
$ws.Range("I132").Value = 1673.1666
$ws.Range("J132").Value = 3150
$ws.Range("K132").Value = 5019.4998
$ws.Range("L132").Value = 9450
$ws.Range("M132").Value = -2489.4998
$ws.Range("N132").Value = -14510

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 92950
$ws.Range("I40").Value = 251312.5
$ws.Range("J40").Value = 2457.1428
$ws.Range("K40").Value = 251312.5
$ws.Range("L40").Value = 2457.1428
$ws.Range("M40").Value = -251176.5
$ws.Range("N40").Value = -2729.1428
$ws.Range("H41").Value = 9137.444
$ws.Range("I41").Value = 4000
$ws.Range("J41").Value = 9779.625
$ws.Range("K41").Value = 4000
$ws.Range("L41").Value = 9779.625
$ws.Range("M41").Value = -3562
$ws.Range("N41").Value = -10655.625
$ws.Range("H122").Value = 2204
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2204
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6612
$ws.Range("M122").ClearContents() | Out-Null
$ws.Range("N122").Value = -11512
$ws.Range("H132").Value = 2933.6316
$ws.Range("I132").Value = 3321.1428
$ws.Range("J132").Value = 1848.6
$ws.Range("K132").Value = 9963.428400000001
$ws.Range("L132").Value = 5545.799999999999
$ws.Range("M132").Value = -7433.428400000001
$ws.Range("N132").Value = -10605.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 25958
$ws.Range("I21").Value = 9815
$ws.Range("J21").Value = 31339
$ws.Range("K21").Value = 9815
$ws.Range("L21").Value = 31339
$ws.Range("M21").Value = -9580
$ws.Range("N21").Value = -31809
$ws.Range("H35").Value = 25958
$ws.Range("I35").Value = 9815
$ws.Range("J35").Value = 31339
$ws.Range("K35").Value = 9815
$ws.Range("L35").Value = 31339
$ws.Range("M35").Value = -9525
$ws.Range("N35").Value = -31919
$ws.Range("H132").Value = 3691.44
$ws.Range("I132").Value = 4478.385
$ws.Range("J132").Value = 2838.9167
$ws.Range("K132").Value = 13435.155
$ws.Range("L132").Value = 8516.750100000001
$ws.Range("M132").Value = -10905.155
$ws.Range("N132").Value = -13576.7501
